$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the old 4-column, 5-row table and shrink it down to a single
# "name" column with 3 data rows (".NET app 01", "Mobile App", "Web App ").
$ws.Cells.Clear()

$ws.Range("A1").Value = "name"
$ws.Range("A2").Value = ".NET app 01"
$ws.Range("A3").Value = "Mobile App"
$ws.Range("A4").Value = "Web App "
